$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Range 1" label in A2 to the new date label.
$ws.Range("A2").Value = "Date 08-06-2025"

# Clear the header row (row 3) text values, keep styles/formatting intact.
$ws.Range("A3:H3").ClearContents()

# Clear the data rows (4-6) values, keep styles/formatting intact.
$ws.Range("A4:H6").ClearContents()

# Remove row 7 entirely (shifts nothing up since it's the last row).
$ws.Rows.Item(7).Delete()

# Move the active selection to H6 to match the new smaller data range.
$ws.Range("H6").Select()
